# Update the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape from coinranking.com (GitHub Actions scheduled update).
# Column A (rank index) is unchanged; only Coin (B), Link (C), Price (D)
# and Volume(1h) (E) values are refreshed for the affected rows. A few
# rows also had their two adjacent entries swap places (e.g. rows 17/18,
# 28/29, 44/45, 49/50). Price cells that look like plain numbers are
# forced to Text format first so Excel keeps the original string
# formatting (e.g. trailing zeros) instead of normalizing them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "67.731.51"
$ws.Range("E2").Value2 = "  -0.71%  "
$ws.Range("D3").Value2 = "3.792.99"
$ws.Range("E3").Value2 = "  +0.46%  "
$ws.Range("E4").Value2 = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "596.00"
$ws.Range("E5").Value2 = "  +0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "167.12"
$ws.Range("E6").Value2 = "  -0.18%  "
$ws.Range("D7").Value2 = "3.792.76"
$ws.Range("E7").Value2 = "  +0.55%  "
$ws.Range("E8").Value2 = "  +0.03%  "
$ws.Range("E9").Value2 = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.160"
$ws.Range("E10").Value2 = "  -0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "6.33"
$ws.Range("E11").Value2 = "  -1.26%  "
$ws.Range("E12").Value2 = "  -0.34%  "
$ws.Range("E13").Value2 = "  -2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "36.05"
$ws.Range("E14").Value2 = "  -0.17%  "
$ws.Range("D15").Value2 = "4.427.13"
$ws.Range("E15").Value2 = "  +0.71%  "
$ws.Range("D16").Value2 = "3.787.52"
$ws.Range("E16").Value2 = "  +0.58%  "
$ws.Range("B17").Value2 = "Chainlink"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "18.48"
$ws.Range("E17").Value2 = "  +3.24%  "
$ws.Range("B18").Value2 = "WrappedBTC"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value2 = "67.693.74"
$ws.Range("E18").Value2 = "  -0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "7.05"
$ws.Range("E19").Value2 = "  +0.64%  "
$ws.Range("E20").Value2 = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "9.99"
$ws.Range("E21").Value2 = "  -7.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "459.32"
$ws.Range("E22").Value2 = "  -1.34%  "
$ws.Range("E23").Value2 = "  -0.17%  "
$ws.Range("E24").Value2 = "  +1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "83.41"
$ws.Range("E25").Value2 = "  -0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "12.08"
$ws.Range("E26").Value2 = "  +1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.11"
$ws.Range("E27").Value2 = "  -3.27%  "
$ws.Range("B28").Value2 = "Dai"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "1.00"
$ws.Range("E28").Value2 = "  -0.10%  "
$ws.Range("B29").Value2 = "RenderToken"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "10.00"
$ws.Range("E29").Value2 = "  -0.94%  "
$ws.Range("E30").Value2 = "  -0.29%  "
$ws.Range("E31").Value2 = "  +4.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "7.21"
$ws.Range("E32").Value2 = "  -1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "29.64"
$ws.Range("E33").Value2 = "  -1.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.999"
$ws.Range("E34").Value2 = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "9.06"
$ws.Range("E35").Value2 = "  -1.00%  "
$ws.Range("D36").Value2 = "3.733.22"
$ws.Range("E36").Value2 = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.100"
$ws.Range("E37").Value2 = "  -0.59%  "
$ws.Range("E38").Value2 = "  -3.18%  "
$ws.Range("E39").Value2 = "  -0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.993"
$ws.Range("E40").Value2 = "  -1.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "5.77"
$ws.Range("E41").Value2 = "  -0.13%  "
$ws.Range("E42").Value2 = "  +0.09%  "
$ws.Range("E43").Value2 = "  +0.03%  "
$ws.Range("B44").Value2 = "OKB"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "48.02"
$ws.Range("E44").Value2 = "  +2.02%  "
$ws.Range("B45").Value2 = "Arweave"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "43.94"
$ws.Range("E45").Value2 = "  -0.84%  "
$ws.Range("E46").Value2 = "  -0.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "149.69"
$ws.Range("E47").Value2 = "  +2.98%  "
$ws.Range("E48").Value2 = "  -1.55%  "
$ws.Range("B49").Value2 = "EnergySwap"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "26.85"
$ws.Range("E49").Value2 = "  +7.05%  "
$ws.Range("B50").Value2 = "Bittensor"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "390.14"
$ws.Range("E50").Value2 = "  +0.42%  "
$ws.Range("E51").Value2 = "  -4.66%  "
